# Automatic update of files.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump the "Förändrad" (column C) date for every existing data row
#    (rows 2..32) from 2024-11-15 (45611) to 2024-11-16 (45612).
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 3).Value = 45612
}

# 2. Give row 32 an explicit row height (matches the OOXML diff which adds
#    ht="15" customHeight="1" to <row r="32">).
$ws.Rows.Item(32).RowHeight = 15

# 3. Append the new row 33 with the new case data.
$row = 33
$ws.Cells.Item($row, 1).Value = "A 52574-2024"   # A: Beteckning
$ws.Cells.Item($row, 2).Value = 45609            # B: Datum
$ws.Cells.Item($row, 3).Value = 45612            # C: Förändrad
$ws.Cells.Item($row, 4).Value = "OKÄNT"          # D: Län
$ws.Cells.Item($row, 5).Value = "OKÄNT"          # E: Kommun
$ws.Cells.Item($row, 7).Value = 0.2              # G: Area (ha)
$ws.Cells.Item($row, 8).Value = 0                # H: Fridlysta
$ws.Cells.Item($row, 9).Value = 0                # I: Signalarter
$ws.Cells.Item($row, 10).Value = 0               # J: NT
$ws.Cells.Item($row, 11).Value = 0               # K: VU
$ws.Cells.Item($row, 12).Value = 0               # L: EN
$ws.Cells.Item($row, 13).Value = 0               # M: CR
$ws.Cells.Item($row, 14).Value = 0               # N: RE
$ws.Cells.Item($row, 15).Value = 0               # O: Rödlistade
$ws.Cells.Item($row, 16).Value = 0               # P: Hotade
$ws.Cells.Item($row, 17).Value = 0               # Q: Alla arter

# Copy the date number-format style from the row above for B/C so the new
# row matches the existing date columns' formatting.
$ws.Cells.Item($row, 2).NumberFormat = $ws.Cells.Item($row - 1, 2).NumberFormat
$ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item($row - 1, 3).NumberFormat

# R33 mirrors the other rows' empty-but-styled "Artnamn" cell: same
# wrap-text style as R2..R32, left blank.
$ws.Cells.Item($row, 18).WrapText = $true
